# Chambers/White material_and_methods.xlsx update
# - coring_method for Chambers_et_al_2019 (row 4) set to "push core"
# - compaction_flag for rows 4 & 5 updated from "not specified" to "no obvious compaction"
# - refresh the sheet's active selection to reflect where the editor was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 (coring_method, Chambers_et_al_2019) was blank -> "push core"
$ws.Range("B4").Value = "push core"

# F4/F5 (compaction_flag) "not specified" -> "no obvious compaction"
$ws.Range("F4").Value = "no obvious compaction"
$ws.Range("F5").Value = "no obvious compaction"

# Move the viewport/selection to match the editor's final cursor position
$ws.Range("Q1").Select()
$ws.Range("F13").Select()
